# PANGEAHIVsim Evaluation Sheet - "last version of the report"
#
# Fills in the NUMBER_SEQ (col D) values that were missing for rows 6-11,
# and corrects the GROWTH (col E) flags so every populated row reads 1
# (row 8 was recorded as 0 and rows 9-11 had no GROWTH value at all).
# Also moves the active selection to D12, the cell just below the data
# entry table, matching where the author's cursor ended up when the file
# was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- NUMBER_SEQ (column D) -------------------------------------------------
$ws.Range("D6").Value  = 200
$ws.Range("D7").Value  = 100
$ws.Range("D8").Value  = 100
$ws.Range("D9").Value  = 200
$ws.Range("D10").Value = 300
$ws.Range("D11").Value = 200

# --- GROWTH (column E) ------------------------------------------------------
$ws.Range("E8").Value  = 1
$ws.Range("E9").Value  = 1
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1

# --- cursor / selection -----------------------------------------------------
$ws.Range("D12").Select()
